$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 254, shifting rows 254:361 down to 255:362
$ws.Range("A254").EntireRow.Insert()

# Populate the newly inserted row 254 with the new record's data.
# Values for A,B,C,E,F,G,H,I,K,L,M,N,P,Q,R are carried over from the
# record that used to occupy row 254 (now at row 255); D, J and O hold
# the new data for this record.
$ws.Range("A254").Value = 4
$ws.Range("B254").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C254").Value = "Los Lagos"
$ws.Range("D254").Value = 44825
$ws.Range("E254").Value = 10
$ws.Range("F254").Value = 100112045
$ws.Range("G254").Value = "Zapallo"
$ws.Range("H254").Value = "Paine"
$ws.Range("I254").Value = "1a (guarda)"
$ws.Range("J254").Value = 500
$ws.Range("K254").Value = 600
$ws.Range("L254").Value = 600
$ws.Range("M254").Value = 600
$ws.Range("N254").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O254").Value = "Región de O'Higgins"
$ws.Range("P254").Value = 600
$ws.Range("Q254").Value = 1
$ws.Range("R254").Value = "Hortaliza"
